$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2" = 21.38850176573779
    "B2" = 17.37868957931509
    "C2" = 25.36871322633619
    "A3" = 20.64072820321903
    "B3" = 16.0178969511467
    "C3" = 26.10571993187365
    "A4" = 20.66920240661135
    "B4" = 16.35110799264588
    "C4" = 25.14206235033802
    "A5" = 24.87244596593988
    "B5" = 19.64858264946364
    "C5" = 30.84565315845443
    "A6" = 26.06727152267131
    "B6" = 21.33245820557327
    "C6" = 29.92419330030981
    "A7" = 19.22750204057473
    "B7" = 16.30796156276427
    "C7" = 22.39451297710363
    "A8" = 9.890052560229611
    "B8" = 4.746148511182452
    "C8" = 16.03020413814457
    "A9" = 24.44687836609264
    "B9" = 20.24378479393847
    "C9" = 28.62347068908447
    "A10" = 10.17638673681832
    "B10" = 5.031857197501548
    "C10" = 16.35607351172702
    "A11" = 17.91087570579169
    "B11" = 13.20271457375548
    "C11" = 22.40668834637858
    "A12" = 26.89164199088832
    "B12" = 21.89782605212991
    "C12" = 31.10846086657531
    "A13" = 12.78397068966698
    "B13" = 6.932317368398227
    "C13" = 20.92502569671525
    "A14" = 14.85808382223082
    "B14" = 9.632307518776107
    "C14" = 20.26029453142244
    "A15" = 28.77413975488276
    "B15" = 22.82338897267322
    "C15" = 34.74708717883421
    "A16" = 26.1384685591192
    "B16" = 21.25317350696336
    "C16" = 29.91458930443408
    "A17" = 26.99126459120875
    "B17" = 22.4963782731673
    "C17" = 31.2811190237013
    "A18" = 20.60665854452809
    "B18" = 16.33698341913452
    "C18" = 25.03073071792736
    "A19" = 18.68931229132765
    "B19" = 13.23803387561573
    "C19" = 24.22573925620963
    "A20" = 29.56814665555689
    "B20" = 24.37545715621447
    "C20" = 34.32573190242808
    "A21" = 28.89050693472969
    "B21" = 23.93114122656443
    "C21" = 34.07804077971073
    "A22" = 10.74804460589509
    "B22" = 5.250898458035385
    "C22" = 17.55835583889407
    "A23" = 30.99345122581428
    "B23" = 24.15346239898594
    "C23" = 37.29133054817315
    "A24" = 26.1384685591192
    "B24" = 21.25317350696336
    "C24" = 29.91458930443408
    "A25" = 22.55094707931317
    "B25" = 19.33264381390847
    "C25" = 26.32194092216196
    "A26" = 28.85021377908337
    "B26" = 23.61260539834298
    "C26" = 34.034651468668
    "A27" = 12.86914838081217
    "B27" = 6.955230536556061
    "C27" = 21.02652352109294
    "A28" = 20.22633749882602
    "B28" = 16.68710789715158
    "C28" = 23.49071654528738
    "A29" = 32.66265051256413
    "B29" = 26.27604153603249
    "C29" = 39.53480947860753
    "A30" = 25.50739117962885
    "B30" = 21.58851508969947
    "C30" = 30.31352674972676
    "A31" = 14.33504672804971
    "B31" = 9.171508245504533
    "C31" = 19.57510881383944
    "A32" = 11.98886714649602
    "B32" = 7.025345031354191
    "C32" = 17.50896397477052
    "A33" = 14.85980855250412
    "B33" = 9.700676961623293
    "C33" = 20.68343499651575
    "A34" = 23.28659907095988
    "B34" = 18.38033576939813
    "C34" = 28.45723910516685
    "A35" = 27.91451393216179
    "B35" = 23.1775276271396
    "C35" = 32.51117257032084
    "A36" = 19.14187521535924
    "B36" = 16.24790871493871
    "C36" = 22.56157526471678
    "A37" = 30.84854018950216
    "B37" = 24.89954793848969
    "C37" = 36.30000920081681
    "A38" = 15.06993194250352
    "B38" = 10.21714930178505
    "C38" = 21.91996661852801
    "A39" = 27.41040841398803
    "B39" = 22.46689449010432
    "C39" = 31.9405343393089
    "A40" = 24.96961133727889
    "B40" = 20.23237721451717
    "C40" = 30.32219849277783
    "A41" = 18.12944049482694
    "B41" = 15.35839246348332
    "C41" = 21.33164893577558
    "A42" = 18.24211456898124
    "B42" = 12.85226734548205
    "C42" = 23.36691001988499
    "A43" = 25.504640455482
    "B43" = 21.58804911889032
    "C43" = 30.31297569725686
    "A44" = 15.92201849085944
    "B44" = 12.72391477802549
    "C44" = 19.80694731927623
    "A45" = 18.88210395946783
    "B45" = 13.58726727389485
    "C45" = 24.57389336615204
    "A46" = 13.92789657152717
    "B46" = 9.895700210756637
    "C46" = 18.00861677089005
    "A47" = 10.72204594366823
    "B47" = 4.778532245218436
    "C47" = 17.54652109482504
    "A48" = 27.01261573401585
    "B48" = 22.21761752663176
    "C48" = 31.38910385922219
    "A49" = 23.02109739769512
    "B49" = 18.75886044642751
    "C49" = 27.72936049084308
    "A50" = 30.88686727013537
    "B50" = 25.08618065526323
    "C50" = 36.34658014902942
    "A51" = 20.66571531911547
    "B51" = 15.93472118723952
    "C51" = 26.37157351736791
    "A52" = 13.70106460147985
    "B52" = 9.439008503937766
    "C52" = 18.07050135695347
    "A53" = 27.45249729259701
    "B53" = 23.03205505260803
    "C53" = 32.17296359274947
    "A54" = 28.76415308248434
    "B54" = 22.67004976782209
    "C54" = 34.78157471428933
    "A55" = 28.89867935509088
    "B55" = 23.56485911861971
    "C55" = 33.86429737170005
    "A56" = 14.108150434489
    "B56" = 9.286287986688228
    "C56" = 19.70038115335174
    "A57" = 11.99581729003972
    "B57" = 7.02175714922909
    "C57" = 17.50716359696655
    "A58" = 30.80002138178461
    "B58" = 25.10211111247232
    "C58" = 36.32106147232457
    "A59" = 17.16051106517111
    "B59" = 12.79420065015565
    "C59" = 21.4318567490837
    "A60" = 31.10796151835798
    "B60" = 24.22863171642936
    "C60" = 37.40628829274337
    "A61" = 14.93298572776943
    "B61" = 9.479863190295289
    "C61" = 21.40704380969573
    "A62" = 17.91087570579169
    "B62" = 13.20271457375548
    "C62" = 22.40668834637858
    "A63" = 15.49707758940566
    "B63" = 9.716216863684249
    "C63" = 22.46081801506012
    "A64" = 30.99710692234159
    "B64" = 24.16986619413412
    "C64" = 37.17711598761869
    "A65" = 11.39037787020175
    "B65" = 6.97900815168977
    "C65" = 16.74619532110268
    "A66" = 22.83286725183245
    "B66" = 18.28256754449947
    "C66" = 27.28715712519181
    "A67" = 22.02165783175921
    "B67" = 17.97020714102748
    "C67" = 26.49260862285382
    "A68" = 26.99206492754991
    "B68" = 21.60037042269972
    "C68" = 32.1332789757786
    "A69" = 10.90046246616398
    "B69" = 5.372930118387
    "C69" = 17.99892818980382
    "A70" = 24.87181308967304
    "B70" = 19.73722240448686
    "C70" = 30.25865547394634
    "A71" = 19.83276672779148
    "B71" = 16.74676168060025
    "C71" = 22.98989762669703
    "A72" = 22.87182159394246
    "B72" = 19.12947158499937
    "C72" = 26.73778111008175
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}